# Update the workbook per the commit: refresh the "Förändrad" timestamp
# column for all existing records and append three newly reported cases.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) All existing data rows (2-387): column C ("Förändrad") changes
#    from 45172 to 45175.
$ws.Range("C2:C387").Value2 = 45175

# 2) Row 387 gains an explicit row height (15pt, customHeight) - matches
#    the new rows appended below.
$ws.Rows.Item(387).RowHeight = 15

# 3) Append three new report rows (388-390).
$newRows = @(
    @{Row=388; A="A 41359-2023"; B=45174; C=45175; G=5.1},
    @{Row=389; A="A 41353-2023"; B=45174; C=45175; G=2.5},
    @{Row=390; A="A 41360-2023"; B=45174; C=45175; G=1.2}
)

foreach ($r in $newRows) {
    $row = $r.Row

    $ws.Cells.Item($row, 1).Value2 = $r.A

    $bCell = $ws.Cells.Item($row, 2)
    $bCell.Value2 = $r.B
    $bCell.NumberFormat = "YYYY-MM-DD"

    $cCell = $ws.Cells.Item($row, 3)
    $cCell.Value2 = $r.C
    $cCell.NumberFormat = "YYYY-MM-DD"

    $ws.Cells.Item($row, 4).Value2 = "VÄSTERBOTTENS LÄN"
    $ws.Cells.Item($row, 5).Value2 = "MALÅ"

    $ws.Cells.Item($row, 7).Value2 = $r.G

    for ($col = 8; $col -le 17; $col++) {
        $ws.Cells.Item($row, $col).Value2 = 0
    }

    # Column R keeps the sheet's wrap-text style even though it stays blank.
    $ws.Cells.Item($row, 18).WrapText = $true
}

# Rows 388 and 389 carry an explicit row height like row 387; row 390
# (the final row) keeps the default height, matching the source diff.
$ws.Rows.Item(388).RowHeight = 15
$ws.Rows.Item(389).RowHeight = 15
